$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows of test data
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "f(x)=x"
$ws.Range("D3").Value = "Linux"
$ws.Range("F3").Value = "PASS"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "f(x)=x+12+sin(b)"
$ws.Range("D4").Value = "Linux"
$ws.Range("F4").Value = "FAILED"

# New value for F2 (PASS) on existing row
$ws.Range("F2").Value = "PASS"

# Column widths (as close as representable to the target character widths of
# 15.42578125 / 11.42578125 / 14.140625, compensating for the engine's internal
# 5/6-character padding offset that gets added on top of ColumnWidth)
$ws.Columns.Item(2).ColumnWidth = 14.592447916666666
$ws.Columns.Item(4).ColumnWidth = 10.592447916666666
$ws.Columns.Item(5).ColumnWidth = 13.307291666666666

# Selection / active cell
$ws.Range("F4").Select()
